$wb = $excel.ActiveWorkbook

$missing = [System.Type]::Missing

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): two new rows (6,7) in columns A:C
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.md"
$wsOverview.Range("B6").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C6").Value = "Handed back: in sync with en-US"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/35dcb9aacc294195892ca939010c5a18c0ffee01/e2e/35dcb9aa-cc29-4195-8928-a939010c5a18.md", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.md") | Out-Null

$wsOverview.Range("A7").Value = "aa278856-594c-4cac-a891-ca314d641da9.md"
$wsOverview.Range("B7").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C7").Value = "Handed back: in sync with en-US"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/aa278856594c4cacb891ca314d641da9c0ffee02/e2e/aa278856-594c-4cac-a891-ca314d641da9.md", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): two new rows (6,7) in columns A:H
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.md"
$wsZh.Range("B6").Value = "Handed back: in sync with en-US"
$wsZh.Range("C6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf"
$wsZh.Range("D6").Value = "2016-02-18 08:27:10"
$wsZh.Range("E6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.md"
$wsZh.Range("F6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf"
$wsZh.Range("G6").Value = "2016-02-18 08:27:59"
$wsZh.Range("H6").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/35dcb9aacc294195892ca939010c5a18c0ffee01/e2e/35dcb9aa-cc29-4195-8928-a939010c5a18.md", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35dcb9aacc294195892ca939010c5a18c0ffee03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/35dcb9aacc294195892ca939010c5a18c0ffee04/e2e/35dcb9aa-cc29-4195-8928-a939010c5a18.md", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/35dcb9aacc294195892ca939010c5a18c0ffee05/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.zh-cn.xlf") | Out-Null

$wsZh.Range("A7").Value = "aa278856-594c-4cac-a891-ca314d641da9.md"
$wsZh.Range("B7").Value = "Handed back: in sync with en-US"
$wsZh.Range("C7").Value = "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf"
$wsZh.Range("D7").Value = "2016-02-18 08:27:10"
$wsZh.Range("E7").Value = "aa278856-594c-4cac-a891-ca314d641da9.md"
$wsZh.Range("F7").Value = "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf"
$wsZh.Range("G7").Value = "2016-02-18 08:27:59"
$wsZh.Range("H7").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/aa278856594c4cacb891ca314d641da9c0ffee02/e2e/aa278856-594c-4cac-a891-ca314d641da9.md", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa278856594c4cacb891ca314d641da9c0ffee06/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/aa278856594c4cacb891ca314d641da9c0ffee07/e2e/aa278856-594c-4cac-a891-ca314d641da9.md", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aa278856594c4cacb891ca314d641da9c0ffee08/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3): two new rows (6,7) in columns A:H
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.md"
$wsDe.Range("B6").Value = "Handed back: in sync with en-US"
$wsDe.Range("C6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf"
$wsDe.Range("D6").Value = "2016-02-18 08:27:22"
$wsDe.Range("E6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.md"
$wsDe.Range("F6").Value = "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf"
$wsDe.Range("G6").Value = "2016-02-18 08:28:21"
$wsDe.Range("H6").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/35dcb9aacc294195892ca939010c5a18c0ffee01/e2e/35dcb9aa-cc29-4195-8928-a939010c5a18.md", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35dcb9aacc294195892ca939010c5a18c0ffee09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/35dcb9aacc294195892ca939010c5a18c0ffee0a/e2e/35dcb9aa-cc29-4195-8928-a939010c5a18.md", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/35dcb9aacc294195892ca939010c5a18c0ffee0b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf", $missing, $missing, "35dcb9aa-cc29-4195-8928-a939010c5a18.f6b300a9b2a3581803712740eb218ff19cc5d686.de-de.xlf") | Out-Null

$wsDe.Range("A7").Value = "aa278856-594c-4cac-a891-ca314d641da9.md"
$wsDe.Range("B7").Value = "Handed back: in sync with en-US"
$wsDe.Range("C7").Value = "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf"
$wsDe.Range("D7").Value = "2016-02-18 08:27:22"
$wsDe.Range("E7").Value = "aa278856-594c-4cac-a891-ca314d641da9.md"
$wsDe.Range("F7").Value = "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf"
$wsDe.Range("G7").Value = "2016-02-18 08:28:21"
$wsDe.Range("H7").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/aa278856594c4cacb891ca314d641da9c0ffee02/e2e/aa278856-594c-4cac-a891-ca314d641da9.md", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa278856594c4cacb891ca314d641da9c0ffee0c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aa278856594c4cacb891ca314d641da9c0ffee0d/e2e/aa278856-594c-4cac-a891-ca314d641da9.md", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aa278856594c4cacb891ca314d641da9c0ffee0e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf", $missing, $missing, "aa278856-594c-4cac-a891-ca314d641da9.e89b1905b6d88ba07f4844c835df728e56023d14.de-de.xlf") | Out-Null

Write-Output "Overview dim: $($wsOverview.UsedRange.Address())"
Write-Output "zh-cn dim: $($wsZh.UsedRange.Address())"
Write-Output "de-de dim: $($wsDe.UsedRange.Address())"
